$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.774.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.63%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.781.36'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.15%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.55'
$ws.Range('D5').Style = 'Normal'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5131'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.38%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3787'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.47%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07770'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -8.08%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.17'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.54%  '

# Row 11
$ws.Range('E11').Value = '  -2.35%  '

# Row 12
$ws.Range('E12').Value = '  +0.05%  '

# Row 13
$ws.Range('E13').Value = '  -3.77%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.30%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.780.83'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.27%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.153'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.82%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.32'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.62%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001070'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.00%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06558'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.96%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.003'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.07%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.43%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.908'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.98%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.819.12'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.58%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.85%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.236'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.65%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.97'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.21%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.29%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.985.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.17%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.359'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.30%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.19'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.65%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1071'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.24%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.026'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.38%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.625'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.48%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.467'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.75%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07071'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.29%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02315'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.12%  '

# Row 37
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2120'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.98%  '

# Row 38
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.630'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.20%  '

# Row 39
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.998'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.74%  '

# Row 40
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.96%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6072'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.02%  '

# Row 42
$ws.Range('E42').Value = '  -0.03%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.148'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.92%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.321'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.72%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.04'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.19%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5912'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.51%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.709'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.80%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '127.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.45%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.208'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.68%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.892'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.84%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06771'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.99%  '
